# Update gh-pages output data (generated at 456a3b4)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 13197
$ws1.Range("F4").Value = 37
$ws1.Range("F5").Value = 8
$ws1.Range("F6").Value = 107
$ws1.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202409/CZPJLWBD1725345694204.jpeg"
$ws1.Range("F11").Value = 13142
$ws1.Range("F13").Value = 567
$ws1.Range("F14").Value = 8816
$ws1.Range("F15").Value = 7889
$ws1.Range("F21").Value = 6
$ws1.Range("F26").Value = 195
$ws1.Range("F27").Value = 79
$ws1.Range("F28").Value = 346

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 13197
$ws4.Range("F5").Value = 37
$ws4.Range("F6").Value = 8
$ws4.Range("F7").Value = 107
$ws4.Range("I7").Value = "//i1.hdslb.com/bfs/openplatform/202409/CZPJLWBD1725345694204.jpeg"
$ws4.Range("F12").Value = 13142
$ws4.Range("F14").Value = 567
$ws4.Range("F15").Value = 8816
$ws4.Range("F16").Value = 7889
$ws4.Range("F22").Value = 6
$ws4.Range("F29").Value = 195
$ws4.Range("F30").Value = 79
$ws4.Range("F31").Value = 346
